$d = $word.ActiveDocument
$bullet = [char]0x2022

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading and the six bullet
# paragraphs that follow its "Impact" sub-heading, by scanning for the
# known original bullet text (several of these strings are duplicated
# elsewhere in the resume, so paragraph-scoped Range edits are used
# instead of a document-wide Find/Replace).
$targetTexts = @(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%",
    "Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations",
    "Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis",
    "Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality",
    "Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets",
    "Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy"
)

$startIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.Trim()
    $expected = ($bullet + " " + $targetTexts[0]).Trim()
    if ($txt -eq $expected) {
        $startIndex = $i
        break
    }
}

if ($startIndex -eq 0) {
    Write-Output "ERROR: could not locate Key Achievements bullets"
} else {
    # Delete the 4th and 5th bullets (longitudinal data analysis methods,
    # ETL pipelines) outright - highest index first so earlier indices
    # remain valid while we work.
    $d.Paragraphs.Item($startIndex + 4).Range.Delete()
    $d.Paragraphs.Item($startIndex + 3).Range.Delete()

    # Rewrite the remaining four bullets as impact-focused accomplishment
    # statements.
    $d.Paragraphs.Item($startIndex).Range.Text = $bullet + " Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
    $d.Paragraphs.Item($startIndex + 1).Range.Text = $bullet + " `$4.7M savings enabled nonprofit access"
    $d.Paragraphs.Item($startIndex + 2).Range.Text = $bullet + " Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"
    $d.Paragraphs.Item($startIndex + 3).Range.Text = $bullet + " 178% accuracy improvement in racial classification algorithms"

    Write-Output "OK"
}
